# B6-PowerPoint.pptx edit
# 1) Re-style the three tables (slides 14-16) with the new built-in table style.
# 2) Switch the deck's applied theme colours from "Integral / Red Violet"
#    back to the default "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style change -------------------------------------------------
$newStyleId = "{5BAB4EA8-2A00-4D12-9B2E-F3ED88A34C56}"
foreach ($slideIdx in 14,15,16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colour change --------------------------------------------------
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Office Theme colour scheme (RGB values, stored as BGR-packed ints like VBA's .RGB)
$colors.Item(1).RGB  = 0         # Dark 1    - 000000
$colors.Item(2).RGB  = 16777215  # Light 1   - FFFFFF
$colors.Item(3).RGB  = 6968388   # Dark 2    - 44546A
$colors.Item(4).RGB  = 15132391  # Light 2   - E7E6E6
$colors.Item(5).RGB  = 13998939  # Accent 1  - 5B9BD5
$colors.Item(6).RGB  = 3243501   # Accent 2  - ED7D31
$colors.Item(7).RGB  = 10855845  # Accent 3  - A5A5A5
$colors.Item(8).RGB  = 49407     # Accent 4  - FFC000
$colors.Item(9).RGB  = 12874308  # Accent 5  - 4472C4
$colors.Item(10).RGB = 4697456   # Accent 6  - 70AD47
$colors.Item(11).RGB = 12673797  # Hyperlink - 0563C1
$colors.Item(12).RGB = 7491477   # Followed Hyperlink - 954F72
